# إضافة حدث جديد في Card21 by admin at 2025-12-08 08:20:04
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card21")

# Row 15 previously had several "empty" placeholder cells (B:K, M) that were
# never filled in when it was created as the most-recent event row. Now that
# a new event is being appended, those placeholders get the same "nan" text
# used throughout the rest of the table.
$nanCols = @("B","C","D","E","F","G","H","I","J","K","M")
foreach ($c in $nanCols) {
    $ws.Range($c + "15").Value = "nan"
}

# Append the new event as row 16, following the same layout as every other
# row: card number (as text, matching the rest of column A) in A16, the
# placeholder columns B:K and M16 left blank, and the date/correction/
# serviced-by details filled in for the new event.
$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = "21"
$ws.Range("A16").Style = "Normal"

foreach ($c in $nanCols) {
    $ws.Range($c + "16").Style = "Normal"
}

$ws.Range("L16").Value = "6\11\2024"
$ws.Range("N16").Value = "تم عمل مرجعه علي معيار المكنه بسبب مشكله في cv"
$ws.Range("O16").Value = "الخبير"
